$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = "aa"
$ws.Range("J6").Value = "Agree/Accept"
$ws.Range("I18").Value = "sv"
$ws.Range("J18").Value = "Statement-opinion"
$ws.Range("I33").Value = "sv"
$ws.Range("J33").Value = "Statement-opinion"
$ws.Range("I37").Value = "sd"
$ws.Range("J37").Value = "Statement-non-opinion"
$ws.Range("I41").Value = "sd"
$ws.Range("J41").Value = "Statement-non-opinion"
$ws.Range("I43").Value = "aa"
$ws.Range("J43").Value = "Agree/Accept"
$ws.Range("I45").Value = "aa"
$ws.Range("J45").Value = "Agree/Accept"
$ws.Range("I47").Value = "aa"
$ws.Range("J47").Value = "Agree/Accept"
$ws.Range("I52").Value = "sv"
$ws.Range("J52").Value = "Statement-opinion"
$ws.Range("I62").Value = "sv"
$ws.Range("J62").Value = "Statement-opinion"
$ws.Range("I78").Value = "sd"
$ws.Range("J78").Value = "Statement-non-opinion"
$ws.Range("I82").Value = "aa"
$ws.Range("J82").Value = "Agree/Accept"
$ws.Range("I92").Value = "aa"
$ws.Range("J92").Value = "Agree/Accept"
$ws.Range("I94").Value = "aa"
$ws.Range("J94").Value = "Agree/Accept"
$ws.Range("I98").Value = "%"
$ws.Range("J98").Value = "Uninterpretable"
$ws.Range("I101").Value = "sd"
$ws.Range("J101").Value = "Statement-non-opinion"
$ws.Range("I105").Value = "sd"
$ws.Range("J105").Value = "Statement-non-opinion"
$ws.Range("I120").Value = "sd"
$ws.Range("J120").Value = "Statement-non-opinion"
$ws.Range("I130").Value = "b"
$ws.Range("J130").Value = "Acknowledge (Backchannel)"
$ws.Range("I143").Value = "sd"
$ws.Range("J143").Value = "Statement-non-opinion"
$ws.Range("I149").Value = "b"
$ws.Range("J149").Value = "Acknowledge (Backchannel)"
$ws.Range("I155").Value = "sd"
$ws.Range("J155").Value = "Statement-non-opinion"
$ws.Range("I165").Value = "sd"
$ws.Range("J165").Value = "Statement-non-opinion"
$ws.Range("I175").Value = "sv"
$ws.Range("J175").Value = "Statement-opinion"
$ws.Range("I177").Value = "ba"
$ws.Range("J177").Value = "Appreciation"
$ws.Range("I178").Value = "aa"
$ws.Range("J178").Value = "Agree/Accept"
$ws.Range("I199").Value = "ba"
$ws.Range("J199").Value = "Appreciation"
$ws.Range("I201").Value = "sv"
$ws.Range("J201").Value = "Statement-opinion"
$ws.Range("I216").Value = "sd"
$ws.Range("J216").Value = "Statement-non-opinion"
$ws.Range("I219").Value = "ba"
$ws.Range("J219").Value = "Appreciation"
$ws.Range("I225").Value = "sd"
$ws.Range("J225").Value = "Statement-non-opinion"
$ws.Range("I230").Value = "ba"
$ws.Range("J230").Value = "Appreciation"
$ws.Range("I231").Value = "aa"
$ws.Range("J231").Value = "Agree/Accept"
$ws.Range("I235").Value = "%"
$ws.Range("J235").Value = "Uninterpretable"
$ws.Range("I245").Value = "sd"
$ws.Range("J245").Value = "Statement-non-opinion"
$ws.Range("I248").Value = "sv"
$ws.Range("J248").Value = "Statement-opinion"
$ws.Range("I252").Value = "sv"
$ws.Range("J252").Value = "Statement-opinion"
$ws.Range("I253").Value = "sv"
$ws.Range("J253").Value = "Statement-opinion"
$ws.Range("I259").Value = "aa"
$ws.Range("J259").Value = "Agree/Accept"
$ws.Range("I281").Value = "aa"
$ws.Range("J281").Value = "Agree/Accept"
$ws.Range("I282").Value = "qy"
$ws.Range("J282").Value = "Yes-No-Question"
$ws.Range("I291").Value = "aa"
$ws.Range("J291").Value = "Agree/Accept"
